# MASTER: IMPLEMENT SUBJECTIVE EVALUATION
#
# Fills in the previously-blank third rater's score (columns E, J, O, T -
# the "TV" rater for each of the four recommenders RB/SS/TV/Global) for all
# 20 evaluated tracks (rows 3-22). The dependent AVERAGE/STDEV formulas in
# columns F,G,K,L,P,Q,U,V and the per-row aggregates in Z/AC, plus the
# precision percentages in F23/K23/P23/U23, recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = 3
$ws.Range("J3").Value  = 3
$ws.Range("O3").Value  = 3
$ws.Range("T3").Value  = 5

$ws.Range("E4").Value  = 2
$ws.Range("J4").Value  = 3
$ws.Range("O4").Value  = 3
$ws.Range("T4").Value  = 5

$ws.Range("E5").Value  = 5
$ws.Range("J5").Value  = 5
$ws.Range("O5").Value  = 4
$ws.Range("T5").Value  = 4

$ws.Range("E6").Value  = 4
$ws.Range("J6").Value  = 2
$ws.Range("O6").Value  = 5
$ws.Range("T6").Value  = 5

$ws.Range("E7").Value  = 2
$ws.Range("J7").Value  = 2
$ws.Range("O7").Value  = 4
$ws.Range("T7").Value  = 5

$ws.Range("E8").Value  = 1
$ws.Range("J8").Value  = 1
$ws.Range("O8").Value  = 4
$ws.Range("T8").Value  = 1

$ws.Range("E9").Value  = 2
$ws.Range("J9").Value  = 3
$ws.Range("O9").Value  = 4
$ws.Range("T9").Value  = 5

$ws.Range("E10").Value = 3
$ws.Range("J10").Value = 3
$ws.Range("O10").Value = 2
$ws.Range("T10").Value = 1

$ws.Range("E11").Value = 1
$ws.Range("J11").Value = 3
$ws.Range("O11").Value = 1
$ws.Range("T11").Value = 1

$ws.Range("E12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("O12").Value = 3
$ws.Range("T12").Value = 1

$ws.Range("E13").Value = 2
$ws.Range("J13").Value = 2
$ws.Range("O13").Value = 3
$ws.Range("T13").Value = 2

$ws.Range("E14").Value = 2
$ws.Range("J14").Value = 5
$ws.Range("O14").Value = 1
$ws.Range("T14").Value = 1

$ws.Range("E15").Value = 2
$ws.Range("J15").Value = 4
$ws.Range("O15").Value = 4
$ws.Range("T15").Value = 1

$ws.Range("E16").Value = 1
$ws.Range("J16").Value = 4
$ws.Range("O16").Value = 1
$ws.Range("T16").Value = 5

$ws.Range("E17").Value = 3
$ws.Range("J17").Value = 2
$ws.Range("O17").Value = 3
$ws.Range("T17").Value = 5

$ws.Range("E18").Value = 1
$ws.Range("J18").Value = 3
$ws.Range("O18").Value = 4
$ws.Range("T18").Value = 5

$ws.Range("E19").Value = 1
$ws.Range("J19").Value = 4
$ws.Range("O19").Value = 2
$ws.Range("T19").Value = 1

$ws.Range("E20").Value = 1
$ws.Range("J20").Value = 3
$ws.Range("O20").Value = 2
$ws.Range("T20").Value = 1

$ws.Range("E21").Value = 4
$ws.Range("J21").Value = 5
$ws.Range("O21").Value = 4
$ws.Range("T21").Value = 5

$ws.Range("E22").Value = 3
$ws.Range("J22").Value = 4
$ws.Range("O22").Value = 3
$ws.Range("T22").Value = 5

# Match the author's final selection in the saved workbook.
$ws.Range("AC29").Select()
